$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 178.66667
$ws.Range("I2").Value = 178.66667
$ws.Range("K2").Value = 178.66667
$ws.Range("M2").Value = -65.66667000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1237304.4
$ws.Range("I33").Value = 1278497.8
$ws.Range("K33").Value = 1278497.8
$ws.Range("M33").Value = -1278268.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1423.5
$ws.Range("I40").Value = 1318.6666
$ws.Range("K40").Value = 1318.6666
$ws.Range("M40").Value = -1143.6666

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 55574668
$ws.Range("I62").Value = 3400
$ws.Range("J62").Value = 125038750
$ws.Range("K62").Value = 3400
$ws.Range("L62").Value = 125038750
$ws.Range("M62").Value = -2776
$ws.Range("N62").Value = -125039998

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 55574668
$ws.Range("I65").Value = 3400
$ws.Range("J65").Value = 125038750
$ws.Range("K65").Value = 17000
$ws.Range("L65").Value = 625193750
$ws.Range("M65").Value = -13880
$ws.Range("N65").Value = -625199990

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1533.75
$ws.Range("I98").Value = 1115.8235
$ws.Range("J98").Value = 3902
$ws.Range("K98").Value = 1115.8235
$ws.Range("L98").Value = 3902
$ws.Range("M98").Value = 382.1765
$ws.Range("N98").Value = -6898

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1533.75
$ws.Range("I122").Value = 1115.8235
$ws.Range("J122").Value = 3902
$ws.Range("K122").Value = 3347.4705
$ws.Range("L122").Value = 11706
$ws.Range("M122").Value = -897.4704999999999
$ws.Range("N122").Value = -16606

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1751.0618
$ws.Range("I132").Value = 1128.5
$ws.Range("J132").Value = 5007.5386
$ws.Range("K132").Value = 3385.5
$ws.Range("L132").Value = 15022.6158
$ws.Range("M132").Value = -855.5
$ws.Range("N132").Value = -20082.6158

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3001114.2
$ws.Range("I137").Value = 1389992
$ws.Range("J137").Value = 7144000
$ws.Range("K137").Value = 4169976
$ws.Range("L137").Value = 21432000
$ws.Range("M137").Value = -4167426
$ws.Range("N137").Value = -21437100

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1236.3684
$ws.Range("I141").Value = 774.7368
$ws.Range("J141").Value = 2621.2632
$ws.Range("K141").Value = 2324.2104
$ws.Range("L141").Value = 7863.7896
$ws.Range("M141").Value = 2855.7896
$ws.Range("N141").Value = -18223.7896

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 76436.69
$ws.Range("I132").Value = 97396.17
$ws.Range("J132").Value = 3777.1333
$ws.Range("K132").Value = 292188.51
$ws.Range("L132").Value = 11331.3999
$ws.Range("M132").Value = -289658.51
$ws.Range("N132").Value = -16391.3999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 38254.4
$ws.Range("J139").Value = 39171.555
$ws.Range("L139").Value = 39171.555
$ws.Range("N139").Value = -49451.555

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 769.4286
$ws.Range("I64").Value = 738.2
$ws.Range("J64").Value = 786.7778
$ws.Range("K64").Value = 738.2
$ws.Range("L64").Value = 786.7778
$ws.Range("M64").Value = -513.2
$ws.Range("N64").Value = -1236.7778

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H67").Value = 769.4286
$ws.Range("I67").Value = 738.2
$ws.Range("J67").Value = 786.7778
$ws.Range("K67").Value = 738.2
$ws.Range("L67").Value = 786.7778
$ws.Range("M67").Value = 41.79999999999995
$ws.Range("N67").Value = -2346.7778

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 814.3684
$ws.Range("I94").Value = 827.2353000000001
$ws.Range("J94").Value = 705
$ws.Range("K94").Value = 827.2353000000001
$ws.Range("L94").Value = 705
$ws.Range("M94").Value = -376.2353000000001
$ws.Range("N94").Value = -1607

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H104").Value = 30000
$ws.Range("J104").Value = 30000
$ws.Range("L104").Value = 30000
$ws.Range("N104").Value = -36988

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 6708.727
$ws.Range("I107").Value = 11636
$ws.Range("J107").Value = 2602.6667
$ws.Range("K107").Value = 11636
$ws.Range("L107").Value = 2602.6667
$ws.Range("M107").Value = -9716
$ws.Range("N107").Value = -6442.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 66.833336
$ws.Range("I7").Value = 50.25
$ws.Range("K7").Value = 50.25
$ws.Range("M7").Value = 62.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1692.6666
$ws.Range("I16").Value = 1642.8572
$ws.Range("J16").Value = 1762.4
$ws.Range("K16").Value = 1642.8572
$ws.Range("L16").Value = 1762.4
$ws.Range("M16").Value = -1355.8572
$ws.Range("N16").Value = -2336.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4056.25
$ws.Range("I62").Value = 5150
$ws.Range("J62").Value = 3400
$ws.Range("K62").Value = 5150
$ws.Range("L62").Value = 3400
$ws.Range("M62").Value = -4526
$ws.Range("N62").Value = -4648

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 4056.25
$ws.Range("I65").Value = 5150
$ws.Range("J65").Value = 3400
$ws.Range("K65").Value = 25750
$ws.Range("L65").Value = 17000
$ws.Range("M65").Value = -22630
$ws.Range("N65").Value = -23240

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1103.1025
$ws.Range("I107").Value = 1146.8518
$ws.Range("J107").Value = 1004.6667
$ws.Range("K107").Value = 1146.8518
$ws.Range("L107").Value = 1004.6667
$ws.Range("M107").Value = 773.1482000000001
$ws.Range("N107").Value = -4844.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1692.6666
$ws.Range("I113").Value = 1642.8572
$ws.Range("J113").Value = 1762.4
$ws.Range("K113").Value = 1642.8572
$ws.Range("L113").Value = 1762.4
$ws.Range("M113").Value = 527.1428000000001
$ws.Range("N113").Value = -6102.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1819.303
$ws.Range("I132").Value = 1731.2885
$ws.Range("J132").Value = 2146.2144
$ws.Range("K132").Value = 5193.8655
$ws.Range("L132").Value = 6438.6432
$ws.Range("M132").Value = -2663.8655
$ws.Range("N132").Value = -11498.6432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H133").Value = 30141.25
$ws.Range("J133").Value = 30141.25
$ws.Range("L133").Value = 30141.25
$ws.Range("N133").Value = -35201.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 533658.2
$ws.Range("I5").Value = 338.41666
$ws.Range("K5").Value = 1015.24998
$ws.Range("M5").Value = -903.2499799999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 48642.01
$ws.Range("I122").Value = 300.30768
$ws.Range("J122").Value = 57021.24
$ws.Range("K122").Value = 2702.76912
$ws.Range("L122").Value = 513191.16
$ws.Range("M122").Value = -252.7691199999999
$ws.Range("N122").Value = -518091.16

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 533658.2
$ws.Range("I135").Value = 338.41666
$ws.Range("K135").Value = 3045.74994
$ws.Range("M135").Value = -510.7499399999997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 530.79486
$ws.Range("I107").Value = 426.22726
$ws.Range("J107").Value = 666.1177
$ws.Range("K107").Value = 426.22726
$ws.Range("L107").Value = 666.1177
$ws.Range("M107").Value = 1493.77274
$ws.Range("N107").Value = -4506.1177

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1829.5625
$ws.Range("I132").Value = 1598.3889
$ws.Range("J132").Value = 2126.7856
$ws.Range("K132").Value = 4795.1667
$ws.Range("L132").Value = 6380.3568
$ws.Range("M132").Value = -2265.1667
$ws.Range("N132").Value = -11440.3568

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 11437.667
$ws.Range("J134").Value = 11437.667
$ws.Range("L134").Value = 34313.001
$ws.Range("N134").Value = -39383.001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1399.9412
$ws.Range("I46").Value = 1141.3572
$ws.Range("J46").Value = 2606.6667
$ws.Range("K46").Value = 1141.3572
$ws.Range("L46").Value = 2606.6667
$ws.Range("M46").Value = -953.3571999999999
$ws.Range("N46").Value = -2982.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2357
$ws.Range("I61").Value = 1939.3846
$ws.Range("J61").Value = 4166.6665
$ws.Range("K61").Value = 1939.3846
$ws.Range("L61").Value = 4166.6665
$ws.Range("M61").Value = -1737.3846
$ws.Range("N61").Value = -4570.6665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1639
$ws.Range("I68").Value = 1511.92
$ws.Range("K68").Value = 1511.92
$ws.Range("M68").Value = -762.9200000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 1639
$ws.Range("I71").Value = 1511.92
$ws.Range("K71").Value = 7559.6
$ws.Range("M71").Value = -3815.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2357
$ws.Range("I113").Value = 1939.3846
$ws.Range("J113").Value = 4166.6665
$ws.Range("K113").Value = 1939.3846
$ws.Range("L113").Value = 4166.6665
$ws.Range("M113").Value = 230.6153999999999
$ws.Range("N113").Value = -8506.666499999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2818.077
$ws.Range("I122").Value = 2400
$ws.Range("J122").Value = 2905.5813
$ws.Range("K122").Value = 7200
$ws.Range("L122").Value = 8716.743899999999
$ws.Range("M122").Value = -4750
$ws.Range("N122").Value = -13616.7439

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H123").Value = 40330
$ws.Range("J123").Value = 40330
$ws.Range("L123").Value = 40330
$ws.Range("N123").Value = -50130

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2390.9714
$ws.Range("I132").Value = 2064.6667
$ws.Range("J132").Value = 3492.25
$ws.Range("K132").Value = 6194.000100000001
$ws.Range("L132").Value = 10476.75
$ws.Range("M132").Value = -3664.000100000001
$ws.Range("N132").Value = -15536.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 60333
$ws.Range("J133").Value = 60333
$ws.Range("L133").Value = 60333
$ws.Range("N133").Value = -65393

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 41022.64
$ws.Range("I113").Value = 67104.60000000001
$ws.Range("J113").Value = 1899.7
$ws.Range("K113").Value = 201313.8
$ws.Range("L113").Value = 5699.1
$ws.Range("M113").Value = -199143.8
$ws.Range("N113").Value = -10039.1

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 17994.357
$ws.Range("I122").Value = 27364
$ws.Range("K122").Value = 82092
$ws.Range("M122").Value = -79642

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1518.9836
$ws.Range("I132").Value = 1302.909
$ws.Range("J132").Value = 3499.6667
$ws.Range("K132").Value = 3908.727
$ws.Range("L132").Value = 10499.0001
$ws.Range("M132").Value = -1378.727
$ws.Range("N132").Value = -15559.0001
